$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.475.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.832.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.612.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.505.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.534"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("E28").Value = "  +4.97%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.78%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.417.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.536"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.24%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "53.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +21.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0483"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.794"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.744.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.838"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.44%  "
